# Saldo_guide.xlsx - daily refresh of the IClientBalance export
# Renames the sheet to the new export timestamp, bumps every reference
# date in column G from 2024-09-17 to 2024-09-18 (serial 45552 -> 45553),
# and corrects the handful of balance rows whose Vl. Total / Saldo
# Previsto figures were revised in the new export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new export run timestamp
$ws.Name = "IClientBalance-20240918-092848-"

# Bump the "Dt. Referencia" date (column G, rows 2-274) by one day
$ws.Range("G2:G274").Value2 = 45553

# Correct the revised balances (columns E "Saldo Previsto" and H "Vl. Total")
$corrections = @{
    52  = 286.58
    104 = -376.23
    110 = 9.9499999999999993
    129 = 98.32
    143 = 1160.72
    224 = 604.54999999999995
}

foreach ($row in $corrections.Keys) {
    $value = $corrections[$row]
    $ws.Cells.Item($row, 5).Value2 = $value   # column E
    $ws.Cells.Item($row, 8).Value2 = $value   # column H
}
